$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "2021" data column (R) to the table, next to the existing "2020" (Q) column ---

# Header row (row 4): 2021 year label, formatted like the neighboring year cells (Q4)
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("R4").Value = 2021

# Data row (row 5): the 2021 value, formatted like the neighboring data cell (Q5)
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("R5").Value = 72

$ws.Application.CutCopyMode = $false

# --- Update the sheet view: scroll back so column A is visible again, and select R1 ---
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("R1").Select()
